# edit.ps1 - applies the "Add files via upload / Updated Info tab to include
# links to help docs" change to minimum_import.xlsx
#
# Summary of the edit:
#  1. Rename sheet "Explanations" -> "Info"
#  2. Sheet "Minimum": drop the (no longer needed) highlight fill on A1,
#     remove the highlight on the Month/Day header cells (now plain/white),
#     and shorten the submitterID header text.
#  3. Sheet "Info": replace the old explanatory table with two short rows
#     that link out to the online docs and the YouTube channel.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the second sheet
# ---------------------------------------------------------------------
$wsMin = $wb.Worksheets.Item("Minimum")
$wsInfo = $wb.Worksheets.Item("Explanations")
$wsInfo.Name = "Info"

# ---------------------------------------------------------------------
# 2. "Minimum" sheet header row tweaks
# ---------------------------------------------------------------------
# A1 loses its highlight fill (back to the plain/default style)
$wsMin.Range("A1").Style = "Normal"

# Month (E1) and Day (F1) no longer carry the "Time" accent color -
# they become a plain white fill instead of inheriting D1's highlight
$wsMin.Range("E1:F1").Interior.ThemeColor = 2
$wsMin.Range("E1:F1").Interior.TintAndShade = 0

# H1 text is shortened - the parenthetical note moved to the Info sheet
$wsMin.Range("H1").Value = "Encounter.submitterID"

# ---------------------------------------------------------------------
# 3. "Info" sheet - replace the whole explanatory table with two rows
#    of help links
# ---------------------------------------------------------------------
$wsInfo.Cells.Clear()

$wsInfo.Range("A1").Value = "For a list of available fields and how to use them see:"
$wsInfo.Range("B1").Value = "https://wildbook.docs.wildme.org/data/bulk-import-beta.html#fields-available"
$wsInfo.Hyperlinks.Add($wsInfo.Range("B1"), "https://wildbook.docs.wildme.org/data/bulk-import-beta.html#fields-available", "fields-available")

$wsInfo.Range("A2").Value = "Find demos and tutorials on our YouTube channel: "
$wsInfo.Range("B2").Value = "https://www.youtube.com/@wildme3451/videos"
$wsInfo.Hyperlinks.Add($wsInfo.Range("B2"), "https://www.youtube.com/@wildme3451/videos")

# Row heights match the new, shorter content
$wsInfo.Rows.Item(1).RowHeight = 15.75
$wsInfo.Rows.Item(2).RowHeight = 16

# Column widths sized to the new text
$wsInfo.Columns.Item(1).ColumnWidth = 39.83
$wsInfo.Columns.Item(2).ColumnWidth = 62

Write-Host "Edit applied"
